$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select column L (as the user did before deleting it) and delete the entire column,
# which shifts columns M and N (PRODUCT CODE, URL LINK) left to L and M.
$ws.Range("L1:L1048576").Select()
$ws.Range("L1").EntireColumn.Delete()
